# Updated CHE_grids model - 2025-08-10 01:24
# Re-assign the "grid_cell" labels in column AG of the "solar" sheet
# (these back the distr_elc_won-CHE_xxxx connection rows in columns R:AG).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$newValues = @{
    7  = "CHE_0"
    8  = "CHE_17"
    9  = "CHE_19"
    10 = "CHE_3"
    11 = "CHE_12"
    12 = "CHE_24"
    13 = "CHE_5"
    14 = "CHE_8"
    15 = "CHE_14"
    16 = "CHE_18"
    17 = "CHE_7"
    20 = "CHE_21"
    21 = "CHE_9"
    22 = "CHE_4"
    23 = "CHE_20"
    24 = "CHE_1"
    25 = "CHE_6"
    26 = "CHE_13"
}

foreach ($row in $newValues.Keys) {
    $ws.Range("AG$row").Value = $newValues[$row]
}
